$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (E / error) values for rows 3-11: the error values were
# recomputed (split into separate sections of the iterative routine).
$ws.Range("D3").Value = 17.6902860206768
$ws.Range("D4").Value = 13.5173641896727
$ws.Range("D5").Value = 3.51949495555342
$ws.Range("D6").Value = 0.685513898070205
$ws.Range("D7").Value = 0.126050414764144
$ws.Range("D8").Value = 0.0229327615412629
$ws.Range("D9").Value = 0.0041641661015461
$ws.Range("D10").Value = 0.0007558700112753
$ws.Range("D11").Value = 0.000137195062706

# Append two further iterations (rows 12 and 13) that resulted from the
# extra iteration steps.
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 38.5667243729146
$ws.Range("C12").Value = [double]"-4.51970664983037e-06"
$ws.Range("D12").Value = [double]"2.49014605913089e-05"

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 38.5667288926213
$ws.Range("C13").Value = [double]"-8.20343057483797e-07"
$ws.Range("D13").Value = [double]"4.51970664983037e-06"
